$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

$ws.Range("D2").Value = "65.391.48"
$ws.Range("D3").Value = "3.332.51"
$ws.Range("E3").Value = "  -4.49%  "
$ws.Range("E4").Value = "  -0.02%  "
Set-TextValue "D5" "573.99"
$ws.Range("E5").Value = "  -1.52%  "
Set-TextValue "D6" "177.52"
$ws.Range("E6").Value = "  +2.72%  "
Set-TextValue "D7" "0.616"
$ws.Range("E7").Value = "  +3.12%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").Value = "3.329.03"
$ws.Range("E9").Value = "  -4.57%  "
Set-TextValue "D10" "0.128"
$ws.Range("E10").Value = "  -2.01%  "
$ws.Range("E11").Value = "  -0.21%  "
Set-TextValue "D12" "0.407"
$ws.Range("E12").Value = "  -0.90%  "
$ws.Range("D13").Value = "3.908.52"
$ws.Range("E13").Value = "  -4.48%  "
$ws.Range("E14").Value = "  +0.17%  "
Set-TextValue "D15" "28.35"
$ws.Range("E15").Value = "  -5.03%  "
$ws.Range("D16").Value = "65.391.99"
$ws.Range("E16").Value = "  -0.96%  "
$ws.Range("E17").Value = "  -1.88%  "
$ws.Range("D18").Value = "3.322.76"
$ws.Range("E18").Value = "  -4.55%  "
Set-TextValue "D19" "5.74"
$ws.Range("E19").Value = "  -3.03%  "
Set-TextValue "D20" "13.37"
$ws.Range("E20").Value = "  -3.90%  "
Set-TextValue "D21" "361.58"
$ws.Range("E21").Value = "  -1.47%  "
$ws.Range("E22").Value = "  -4.15%  "
$ws.Range("E23").Value = "  -0.11%  "
Set-TextValue "D24" "71.03"
$ws.Range("E24").Value = "  -2.49%  "
Set-TextValue "D25" "0.516"
$ws.Range("E25").Value = "  -3.27%  "
Set-TextValue "D26" "0.0000122"
$ws.Range("E26").Value = "  -3.82%  "
$ws.Range("E27").Value = "  -1.13%  "
$ws.Range("E28").Value = "  -1.08%  "
$ws.Range("E29").Value = "  +0.01%  "
$ws.Range("E30").Value = "  -1.59%  "
Set-TextValue "D31" "0.999"
$ws.Range("E31").Value = "  -0.07%  "
Set-TextValue "D32" "5.58"
$ws.Range("E32").Value = "  -3.18%  "
Set-TextValue "D33" "22.85"
$ws.Range("E33").Value = "  -5.20%  "
Set-TextValue "D34" "6.80"
$ws.Range("E34").Value = "  -4.74%  "
Set-TextValue "D35" "1.21"
$ws.Range("E35").Value = "  -6.38%  "
$ws.Range("E36").Value = "  -3.38%  "
Set-TextValue "D37" "159.95"
$ws.Range("E37").Value = "  -0.53%  "
$ws.Range("E38").Value = "  -5.20%  "
Set-TextValue "D39" "27.31"
$ws.Range("E39").Value = "  -8.28%  "
$ws.Range("E40").Value = "  -0.84%  "
$ws.Range("D41").Value = "2.703.22"
$ws.Range("E41").Value = "  -4.42%  "
$ws.Range("E42").Value = "  -3.16%  "
Set-TextValue "D43" "6.20"
$ws.Range("E43").Value = "  -4.36%  "
$ws.Range("E44").Value = "  -4.48%  "
Set-TextValue "D45" "39.82"
$ws.Range("E45").Value = "  -0.67%  "
$ws.Range("E46").Value = "  -2.87%  "
Set-TextValue "D47" "333.15"
$ws.Range("E47").Value = "  +2.85%  "
Set-TextValue "D48" "23.82"
$ws.Range("E48").Value = "  -1.33%  "
Set-TextValue "D49" "0.0278"
$ws.Range("E49").Value = "  -3.35%  "
Set-TextValue "D50" "0.103"
$ws.Range("E50").Value = "  +1.54%  "
$ws.Range("E51").Value = "  -0.02%  "
